$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 635 (pushes old rows 635..762 down to 636..763)
$ws.Rows.Item(635).Insert()

# Populate the newly inserted row 635 with the new weekly record
$ws.Cells.Item(635, 1).Value = 3
$ws.Cells.Item(635, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(635, 3).Value = "Coquimbo"
$ws.Cells.Item(635, 4).Value = 45275
$ws.Cells.Item(635, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(635, 5).Value = 5
$ws.Cells.Item(635, 6).Value = "Fruta"
$ws.Cells.Item(635, 7).Value = 100108
$ws.Cells.Item(635, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(635, 9).Value = 100108002
$ws.Cells.Item(635, 10).Value = "Mango"
$ws.Cells.Item(635, 11).Value = "Sin especificar"
$ws.Cells.Item(635, 12).Value = "Primera"
$ws.Cells.Item(635, 13).Value = 228
$ws.Cells.Item(635, 14).Value = 10000
$ws.Cells.Item(635, 15).Value = 10000
$ws.Cells.Item(635, 16).Value = 10000
$ws.Cells.Item(635, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(635, 18).Value = "Perú"
$ws.Cells.Item(635, 19).Value = 2500
$ws.Cells.Item(635, 20).Value = 4
